$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Status" column (E) from "COMPLETED" to "success" for all data rows
$ws.Range("E2:E12").Value = "success"

# Clear the "Reason" column (I) for all data rows - the
# "HashKeyGenerated Successfully" note is no longer applicable
$ws.Range("I2:I12").ClearContents()

# Update the active cell selection to reflect where the user left off editing
$ws.Range("H8").Select()
